{"js": "// Remove the \"Scenario: ...\" text from the document's last paragraph\n// while keeping the paragraph (and its bookmark) intact, per:\n//   \"Tolto gli scenari dagli UC\" (Removed the scenarios from the Use Cases)\n//\n// The target text is split across two runs in the source paragraph, so we\n// locate it via Body.search() (which matches across run boundaries) and\n// delete just that text range \u2014 leaving the paragraph's bookmarkStart/\n// bookmarkEnd (and the paragraph mark itself) untouched.\n\nconst searchText =\n  \"Scenario: Un cliente decide di voler visualizzare il proprio carrello. \" +\n  \"Il Sistema fa visualizzare il carrello al Cliente con i prodotti in caso \" +\n  \"ci siano altrimenti mostra un carrello vuoto.\";\n\nconst body = context.document.body;\nconst results = body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the \"Scenario: ...\" text from the document's last paragraph while\n# keeping the paragraph (and its bookmark) intact, per:\n#   \"Tolto gli scenari dagli UC\" (Removed the scenarios from the Use Cases)\n#\n# The target text is split across two runs in the source paragraph, so we\n# locate it with Find.Execute (which matches across run boundaries) and\n# delete just that matched range - leaving the paragraph's bookmarkStart/\n# bookmarkEnd (and the paragraph mark itself) untouched.\n\n$d = $word.ActiveDocument\n\n$searchText = \"Scenario: Un cliente decide di voler visualizzare il proprio carrello. Il Sistema fa visualizzare il carrello al Cliente con i prodotti in caso ci siano altrimenti mostra un carrello vuoto.\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n\n$found = $rng.Find.Execute($searchText)\n\nif ($found) {\n    $rng.Delete()\n}\n"}
